$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in I1 (bold style like H1) and new value in I2
$ws.Range("I1").Value = "pixel_size_mm"
$ws.Range("I1").Font.Bold = $true

$ws.Range("I2").Value = 1.818

# Move the active selection like the saved file shows
[void]$ws.Range("K4").Select()

# Touch page setup so a pageSetup element is emitted
$ws.PageSetup.Orientation = 1
